$wb = $excel.ActiveWorkbook

# --- Sheet "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("E2").Value = 2.7
$ws1.Range("C3").Value = 6
$ws1.Range("E3").Value = 16.2
$ws1.Range("C4").Value = 25
$ws1.Range("E4").Value = 67.59999999999999
$ws1.Range("B5").Value = 4
$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 40
$ws1.Range("E5").Value = 8.1
$ws1.Range("B6").Value = 2
$ws1.Range("D6").Value = 20
$ws1.Range("E6").Value = 5.4
$ws1.Range("C7").Value = 99

# --- Sheet "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("C4").Value = 37

# --- Sheet "Interannual update - High Pri" ---
$ws5 = $wb.Worksheets.Item("Interannual update - High Pri")
$ws5.Range("B2").Value = 60
$ws5.Range("C2").Value = 58.3
$ws5.Range("D2").Value = 60
$ws5.Range("E2").Value = 80
$ws5.Range("B3").Value = 43
$ws5.Range("C3").Value = 41.7
$ws5.Range("D3").Value = 15
$ws5.Range("E3").Value = 20
